$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set H7 to "Move" (PlayOnEnd value for the "Move" animator state row)
$ws.Range("H7").Value = "Move"

# Update the active selection to H7 to match the saved workbook state
$ws.Range("H7").Select()
